$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price/Volume data range so that
# numeric-looking strings (e.g. "374.84") are stored as text, matching
# the source workbook which keeps these as inline/shared strings.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '51.514.20'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.929.93'
$ws.Range('E3').Value = '  -2.56%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '374.84'
$ws.Range('E5').Value = '  +5.71%  '
$ws.Range('D6').Value = '103.07'
$ws.Range('E6').Value = '  -3.76%  '
$ws.Range('D7').Value = '0.543'
$ws.Range('E7').Value = '  -2.71%  '
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').Value = '0.586'
$ws.Range('E9').Value = '  -4.39%  '
$ws.Range('D10').Value = '37.02'
$ws.Range('E10').Value = '  -2.94%  '
$ws.Range('E11').Value = '  -0.53%  '
$ws.Range('E12').Value = '  -2.49%  '
$ws.Range('D13').Value = '18.35'
$ws.Range('E13').Value = '  -3.62%  '
$ws.Range('D14').Value = '3.388.93'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '7.36'
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = '2.924.01'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('D17').Value = '0.931'
$ws.Range('E17').Value = '  -8.42%  '
$ws.Range('D18').Value = '51.427.72'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').Value = '3.44'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').Value = '7.33'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').Value = '12.97'
$ws.Range('E21').Value = '  -4.44%  '
$ws.Range('D22').Value = '0.0₃0946'
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('D23').Value = '68.35'
$ws.Range('E23').Value = '  -1.23%  '
$ws.Range('D24').Value = '262.14'
$ws.Range('E24').Value = '  -0.70%  '
$ws.Range('D25').Value = '2.74'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('E26').Value = '  -5.31%  '
$ws.Range('E27').Value = '  -4.82%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = '25.75'
$ws.Range('E29').Value = '  -4.79%  '
$ws.Range('D30').Value = '7.31'
$ws.Range('E30').Value = '  -1.16%  '
$ws.Range('D31').Value = '6.96'
$ws.Range('E31').Value = '  +8.78%  '
$ws.Range('D32').Value = '0.102'
$ws.Range('E32').Value = '  -4.85%  '
$ws.Range('D33').Value = '9.81'
$ws.Range('E33').Value = '  -3.71%  '
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('D35').Value = '51.05'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('D36').Value = '34.07'
$ws.Range('E36').Value = '  -5.59%  '
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '0.0422'
$ws.Range('E38').Value = '  -3.50%  '
$ws.Range('D39').Value = '3.01'
$ws.Range('E39').Value = '  -9.95%  '
$ws.Range('E40').Value = '  -3.65%  '
$ws.Range('D41').Value = '2.56'
$ws.Range('E41').Value = '  -9.94%  '
$ws.Range('D42').Value = '1.82'
$ws.Range('E42').Value = '  -7.45%  '
$ws.Range('E43').Value = '  -2.39%  '
$ws.Range('D44').Value = '122.03'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('D45').Value = '21.84'
$ws.Range('E45').Value = '  -4.87%  '
$ws.Range('D46').Value = '2.05'
$ws.Range('E46').Value = '  -3.58%  '
$ws.Range('D47').Value = '0.272'
$ws.Range('E47').Value = '  +12.18%  '
$ws.Range('D48').Value = '2.023.76'
$ws.Range('E48').Value = '  -4.69%  '
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').Value = '3.17'
$ws.Range('E50').Value = '  -4.88%  '
$ws.Range('D51').Value = '3.206.70'
$ws.Range('E51').Value = '  -2.92%  '

# Restore the original (default) cell formatting now that the text
# values are committed, so no stray number-format style lingers.
$dataRange.ClearFormats()

Write-Host "Updated cryptos list"
